$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "gcell4"
$ws.Range("B8").Value = "ucell1"

$ws.Range("E16").Select()
